$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the existing "email" hyperlink column from E to D for the rows
#    that already had one, freeing up E for a brand-new "Phone Number"
#    column. Capture the text first, then clear the old cell (value +
#    formatting) before writing it into D.
# ---------------------------------------------------------------------------
$emailRows = 2, 3, 5, 12, 13, 14, 15
foreach ($r in $emailRows) {
    $addr = "E" + $r
    $dst = "D" + $r
    $ws.Range($dst).Value = $ws.Range($addr).Value()
    $ws.Range($addr).Clear()
}

# Drop the old hyperlink definitions (still pointing at column E) so we can
# recreate them against column D in the same original order - this keeps
# the r:id numbering (rId1..rId7) identical to before.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:surin.bachan@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:vikashboodoosingh@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:aasifedoo@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:ravimaharaj_tt@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:bryanjangeesingh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:gtatrinidad@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:jared.hamid@gmail.com")

foreach ($r in $emailRows) {
    $ws.Range("D" + $r).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 2. New header row labels for columns D and E.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Phone Number"

# ---------------------------------------------------------------------------
# 3. Surin Bachan (row 14) also gets a phone number in the new column E.
# ---------------------------------------------------------------------------
$ws.Range("E14").Value = "756-7332"

# ---------------------------------------------------------------------------
# 4. Two new signups appended as rows 21 and 22.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Amar Harripersad"
$ws.Range("B21").Value = "6BS1"
$ws.Range("C21").Value = "Participant"
$ws.Hyperlinks.Add($ws.Range("D21"), "mailto:amarharr@gmail.com")
$ws.Range("D21").Style = "Hyperlink"
$ws.Range("E21").Value = 7046465

$ws.Range("A22").Value = "Sean-Michael Gopaul"
$ws.Range("B22").Value = "4P"
$ws.Range("C22").Value = "Participant"

# ---------------------------------------------------------------------------
# 5. Column widths for the (now populated) D and E columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28.1
$ws.Columns.Item(5).ColumnWidth = 13.1

# ---------------------------------------------------------------------------
# 6. Selection left where the user's cursor ended up after the edit.
# ---------------------------------------------------------------------------
$ws.Range("C23").Select()
